$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record arrived; insert it at row 49, pushing the existing
# rows 49-51 down to 50-52 (their content is preserved unchanged).
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new record's data.
$ws.Cells.Item(49, 1).Value = 9
$ws.Cells.Item(49, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(49, 3).Value = "Metropolitana"
$ws.Cells.Item(49, 4).Value2 = 44615
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(49, 5).Value = 13
$ws.Cells.Item(49, 6).Value = 100112029
$ws.Cells.Item(49, 7).Value = "Orégano"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 16
$ws.Cells.Item(49, 11).Value = 10000
$ws.Cells.Item(49, 12).Value = 10000
$ws.Cells.Item(49, 13).Value = 10000
$ws.Cells.Item(49, 14).Value = "$/docena de atados"
$ws.Cells.Item(49, 15).Value = "Región Metropolitana"
$ws.Cells.Item(49, 16).Value = 3333
$ws.Cells.Item(49, 17).Value = 3
$ws.Cells.Item(49, 18).Value = "Hortaliza"
